$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2243.2307
$ws.Range("J19").Value = 1123.6666
$ws.Range("L19").Value = 1123.6666
$ws.Range("N19").Value = -1473.6666
$ws.Range("H32").Value = 3515.2
$ws.Range("I32").Value = 1253
$ws.Range("J32").Value = 4484.7144
$ws.Range("K32").Value = 1253
$ws.Range("L32").Value = 4484.7144
$ws.Range("M32").Value = -927
$ws.Range("N32").Value = -5136.7144
$ws.Range("H43").Value = 876
$ws.Range("J43").Value = 876
$ws.Range("L43").Value = 876
$ws.Range("N43").Value = -1014
$ws.Range("H62").Value = 1940.1
$ws.Range("I62").Value = 1629.6666
$ws.Range("K62").Value = 1629.6666
$ws.Range("M62").Value = -1005.6666
$ws.Range("H65").Value = 1940.1
$ws.Range("I65").Value = 1629.6666
$ws.Range("K65").Value = 8148.333000000001
$ws.Range("M65").Value = -5028.333000000001
$ws.Range("H76").Value = 5020.4
$ws.Range("I76").Value = 3775.5
$ws.Range("K76").Value = 3775.5
$ws.Range("M76").Value = -3460.5
$ws.Range("H79").Value = 5020.4
$ws.Range("I79").Value = 3775.5
$ws.Range("K79").Value = 3775.5
$ws.Range("M79").Value = -2683.5
$ws.Range("H116").Value = 4494.75
$ws.Range("I116").Value = 4165.25
$ws.Range("J116").Value = 4824.25
$ws.Range("K116").Value = 4165.25
$ws.Range("L116").Value = 4824.25
$ws.Range("M116").Value = -723.25
$ws.Range("N116").Value = -11708.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J32").Value = 778669
$ws.Range("L32").Value = 778669
$ws.Range("N32").Value = -779243
$ws.Range("H61").Value = 1844.3889
$ws.Range("I61").Value = 1844.3889
$ws.Range("K61").Value = 1844.3889
$ws.Range("M61").Value = -1632.3889
$ws.Range("H74").Value = 4918.3335
$ws.Range("I74").Value = 4918.3335
$ws.Range("K74").Value = 4918.3335
$ws.Range("M74").Value = -4044.3335
$ws.Range("H77").Value = 4918.3335
$ws.Range("I77").Value = 4918.3335
$ws.Range("K77").Value = 24591.6675
$ws.Range("M77").Value = -20223.6675
$ws.Range("H122").Value = 1532.8572
$ws.Range("I122").Value = 1232.8334
$ws.Range("K122").Value = 3698.5002
$ws.Range("M122").Value = -1248.5002
$ws.Range("H136").Value = 1844.3889
$ws.Range("I136").Value = 1844.3889
$ws.Range("K136").Value = 5533.1667
$ws.Range("M136").Value = -2983.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 596.625
$ws.Range("I64").Value = 945
$ws.Range("J64").Value = 480.5
$ws.Range("K64").Value = 945
$ws.Range("L64").Value = 480.5
$ws.Range("M64").Value = -720
$ws.Range("N64").Value = -930.5
$ws.Range("H67").Value = 596.625
$ws.Range("I67").Value = 945
$ws.Range("J67").Value = 480.5
$ws.Range("K67").Value = 945
$ws.Range("L67").Value = 480.5
$ws.Range("M67").Value = -165
$ws.Range("N67").Value = -2040.5
$ws.Range("H86").Value = 1394.2858
$ws.Range("I86").Value = 1366.8334
$ws.Range("J86").Value = 1414.875
$ws.Range("K86").Value = 1366.8334
$ws.Range("L86").Value = 1414.875
$ws.Range("M86").Value = -243.8334
$ws.Range("N86").Value = -3660.875
$ws.Range("H89").Value = 1394.2858
$ws.Range("I89").Value = 1366.8334
$ws.Range("J89").Value = 1414.875
$ws.Range("K89").Value = 6834.166999999999
$ws.Range("L89").Value = 7074.375
$ws.Range("M89").Value = -1218.166999999999
$ws.Range("N89").Value = -18306.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 448
$ws.Range("I22").Value = 430.66666
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 430.66666
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -80.66665999999998
$ws.Range("N22").Value = -1200
$ws.Range("H86").Value = 11793.23
$ws.Range("I86").Value = 11868.375
$ws.Range("J86").Value = 11673
$ws.Range("K86").Value = 11868.375
$ws.Range("L86").Value = 11673
$ws.Range("M86").Value = -10745.375
$ws.Range("N86").Value = -13919
$ws.Range("H89").Value = 11793.23
$ws.Range("I89").Value = 11868.375
$ws.Range("J89").Value = 11673
$ws.Range("K89").Value = 59341.875
$ws.Range("L89").Value = 58365
$ws.Range("M89").Value = -53725.875
$ws.Range("N89").Value = -69597
$ws.Range("H124").Value = 90000
$ws.Range("J124").Value = 90000
$ws.Range("L124").Value = 90000
$ws.Range("N124").Value = -94910
$ws.Range("H134").Value = 2922.5557
$ws.Range("I134").Value = 2854
$ws.Range("J134").Value = 3471
$ws.Range("K134").Value = 8562
$ws.Range("L134").Value = 10413
$ws.Range("M134").Value = -6027
$ws.Range("N134").Value = -15483

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10000747
$ws.Range("I4").Value = 10000747
$ws.Range("K4").Value = 30002241
$ws.Range("M4").Value = -30002129
$ws.Range("H40").Value = 74.44444
$ws.Range("I40").Value = 75.71429000000001
$ws.Range("J40").Value = 70
$ws.Range("K40").Value = 302.85716
$ws.Range("L40").Value = 280
$ws.Range("M40").Value = -233.85716
$ws.Range("N40").Value = -418
$ws.Range("H132").Value = 2330.7
$ws.Range("J132").Value = 1832.5714
$ws.Range("L132").Value = 16493.1426
$ws.Range("N132").Value = -21553.1426
$ws.Range("H134").Value = 2041
$ws.Range("I134").Value = 2041
$ws.Range("K134").Value = 6123
$ws.Range("M134").Value = -1053
$ws.Range("H139").Value = 2582
$ws.Range("I139").Value = 2298.6
$ws.Range("J139").Value = 3999
$ws.Range("K139").Value = 6895.799999999999
$ws.Range("L139").Value = 11997
$ws.Range("M139").Value = -1755.799999999999
$ws.Range("N139").Value = -22277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1496
$ws.Range("I122").Value = 1496
$ws.Range("K122").Value = 4488
$ws.Range("M122").Value = -2038
$ws.Range("H126").Value = 2979.8
$ws.Range("I126").Value = 2979.8
$ws.Range("K126").Value = 8939.400000000001
$ws.Range("M126").Value = -6469.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 24910.6
$ws.Range("J74").Value = 25499.5
$ws.Range("L74").Value = 25499.5
$ws.Range("N74").Value = -27371.5
$ws.Range("H75").Value = 74227.25
$ws.Range("J75").Value = 76501
$ws.Range("L75").Value = 76501
$ws.Range("N75").Value = -78373
$ws.Range("H77").Value = 24910.6
$ws.Range("J77").Value = 25499.5
$ws.Range("L77").Value = 76498.5
$ws.Range("N77").Value = -85858.5
$ws.Range("H78").Value = 74227.25
$ws.Range("J78").Value = 76501
$ws.Range("L78").Value = 229503
$ws.Range("N78").Value = -238863
$ws.Range("H136").Value = 2249.6428
$ws.Range("I136").Value = 2291.68
$ws.Range("K136").Value = 6875.039999999999
$ws.Range("M136").Value = -4325.039999999999
